# Adds summary sheets: Total Strut, Total Allthread Cuts, Total Assemblies, Print Me
$wb = $excel.ActiveWorkbook

# --- Create the 4 new worksheets, in order, after the existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsStrut = $wb.Worksheets.Add($null, $lastSheet)
$wsStrut.Name = "Total Strut"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAllthread = $wb.Worksheets.Add($null, $lastSheet)
$wsAllthread.Name = "Total Allthread Cuts"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAssemblies = $wb.Worksheets.Add($null, $lastSheet)
$wsAssemblies.Name = "Total Assemblies"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPrint = $wb.Worksheets.Add($null, $lastSheet)
$wsPrint.Name = "Print Me"

# --- Total Strut ---
$wsStrut.Range("A1").Value = 'Total Strut'
$wsStrut.Range("A1:B2").Merge()
$wsStrut.Range("A1").VerticalAlignment = -4108
$wsStrut.Range("A1").HorizontalAlignment = -4108
$wsStrut.Range("A3").Value = "Strut Type"
$wsStrut.Range("B3").Value = "Quantity"
$wsStrut.Range("A4").Value = 'BACK TO BACK SHALLOW: 14"'
$wsStrut.Range("B4").Value = 1
$wsStrut.Range("A5").Value = 'BACK TO BACK SHALLOW: 22"'
$wsStrut.Range("B5").Value = 1
$wsStrut.Range("A6").Value = 'DEEP STRUT: 12 1/2"'
$wsStrut.Range("B6").Value = 15
$wsStrut.Range("A7").Value = 'DEEP STRUT: 12"'
$wsStrut.Range("B7").Value = 1
$wsStrut.Range("A8").Value = 'DEEP STRUT: 13"'
$wsStrut.Range("B8").Value = 15
$wsStrut.Range("A9").Value = 'DEEP STRUT: 14"'
$wsStrut.Range("B9").Value = 11
$wsStrut.Range("A10").Value = 'DEEP STRUT: 20"'
$wsStrut.Range("B10").Value = 5
$wsStrut.Range("A11").Value = 'DEEP STRUT: 22"'
$wsStrut.Range("B11").Value = 1
$wsStrut.Range("A12").Value = 'DEEP STRUT: 6"'
$wsStrut.Range("B12").Value = 1

# --- Total Allthread Cuts ---
$wsAllthread.Range("A1").Value = 'Total Allthread'
$wsAllthread.Range("A1:B2").Merge()
$wsAllthread.Range("A1").VerticalAlignment = -4108
$wsAllthread.Range("A1").HorizontalAlignment = -4108
$wsAllthread.Range("A3").Value = "Allthread Length"
$wsAllthread.Range("B3").Value = "Quantity"
$wsAllthread.Range("A4").Value = '17 1/2"'
$wsAllthread.Range("B4").Value = 8
$wsAllthread.Range("A5").Value = '22"'
$wsAllthread.Range("B5").Value = 6
$wsAllthread.Range("A6").Value = '23"'
$wsAllthread.Range("B6").Value = 68
$wsAllthread.Range("A7").Value = '25"'
$wsAllthread.Range("B7").Value = 6
$wsAllthread.Range("A8").Value = '27"'
$wsAllthread.Range("B8").Value = 2
$wsAllthread.Range("A9").Value = '28 1/2"'
$wsAllthread.Range("B9").Value = 2
$wsAllthread.Range("A10").Value = '42 1/2"'
$wsAllthread.Range("B10").Value = 8
$wsAllthread.Range("A11").Value = '8"'
$wsAllthread.Range("B11").Value = 2
$wsAllthread.Range("A12").Value = 'Total allthread length = 2453.0'

# --- Total Assemblies ---
$wsAssemblies.Range("A1").Value = 'Assembly Name and Quantity'
$wsAssemblies.Range("A1:B2").Merge()
$wsAssemblies.Range("A1").VerticalAlignment = -4108
$wsAssemblies.Range("A1").HorizontalAlignment = -4108
$wsAssemblies.Range("A3").Value = "Assembly Name"
$wsAssemblies.Range("B3").Value = "Quantity"
$wsAssemblies.Range("A4").Value = 'CR-01'
$wsAssemblies.Range("B4").Value = 1
$wsAssemblies.Range("A5").Value = 'CR-02'
$wsAssemblies.Range("B5").Value = 1
$wsAssemblies.Range("A6").Value = 'CT-03'
$wsAssemblies.Range("B6").Value = 15
$wsAssemblies.Range("A7").Value = 'CT-04'
$wsAssemblies.Range("B7").Value = 15
$wsAssemblies.Range("A8").Value = 'CR-05'
$wsAssemblies.Range("B8").Value = 1
$wsAssemblies.Range("A9").Value = 'CT-06'
$wsAssemblies.Range("B9").Value = 4
$wsAssemblies.Range("A10").Value = 'CT-07'
$wsAssemblies.Range("B10").Value = 3
$wsAssemblies.Range("A11").Value = 'CT-08'
$wsAssemblies.Range("B11").Value = 1
$wsAssemblies.Range("A12").Value = 'CT-09'
$wsAssemblies.Range("B12").Value = 1
$wsAssemblies.Range("A13").Value = 'CT-10'
$wsAssemblies.Range("B13").Value = 2
$wsAssemblies.Range("A14").Value = 'CT-11'
$wsAssemblies.Range("B14").Value = 1
$wsAssemblies.Range("A15").Value = 'CT-12'
$wsAssemblies.Range("B15").Value = 4
$wsAssemblies.Range("A16").Value = 'CT-13'
$wsAssemblies.Range("B16").Value = 2

# --- Print Me ---
$wsPrint.Range("A1").Value = "PRINT_ME"
$wsPrint.Range("A2").Value = 'LEVEL 2 AREA A Tag: CR-01                                                           TOU: 9'' - 4 9/16"                                                                    DEEP STRUT: 6"                                                               Allthread Length: 22"'
$wsPrint.Range("A3").Value = 'LEVEL 2 AREA A Tag: CR-02                                                           TOU: 9'' - 4 9/16"                                                                    DEEP STRUT: 12"                                                               Allthread Length: 22"'
$wsPrint.Range("A4").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A5").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A6").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A7").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A8").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A9").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A10").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A11").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A12").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A13").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A14").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A15").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A16").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A17").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A18").Value = 'LEVEL 2 AREA A Tag: CT-03                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 12 1/2"                                                               Allthread Length: 23"'
$wsPrint.Range("A19").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A20").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A21").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A22").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A23").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A24").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A25").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A26").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A27").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A28").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A29").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A30").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A31").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A32").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A33").Value = 'LEVEL 2 AREA A Tag: CT-04                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 13"                                                               Allthread Length: 23"'
$wsPrint.Range("A34").Value = 'LEVEL 2 AREA A Tag: CR-05                                                           TOU: 9'' - 4 9/16"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 22"'
$wsPrint.Range("A35").Value = 'LEVEL 2 AREA A Tag: CT-06                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 23"'
$wsPrint.Range("A36").Value = 'LEVEL 2 AREA A Tag: CT-06                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 23"'
$wsPrint.Range("A37").Value = 'LEVEL 2 AREA A Tag: CT-06                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 23"'
$wsPrint.Range("A38").Value = 'LEVEL 2 AREA A Tag: CT-06                                                           TOU: 9'' - 3 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 23"'
$wsPrint.Range("A39").Value = 'LEVEL 2 AREA A Tag: CT-07                                                           TOU: 9'' - 1 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 25"'
$wsPrint.Range("A40").Value = 'LEVEL 2 AREA A Tag: CT-07                                                           TOU: 9'' - 1 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 25"'
$wsPrint.Range("A41").Value = 'LEVEL 2 AREA A Tag: CT-07                                                           TOU: 9'' - 1 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 25"'
$wsPrint.Range("A42").Value = 'LEVEL 2 AREA A Tag: CT-08                                                           TOU: 8'' - 11 1/2"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 27"'
$wsPrint.Range("A43").Value = 'LEVEL 2 AREA A Tag: CT-09                                                           TOU: 8'' - 10"                                                                    BACK TO BACK SHALLOW: 14"                                                               Allthread Length: 28 1/2"'
$wsPrint.Range("A44").Value = 'LEVEL 2 AREA A Tag: CT-10                                                           TOU: 7'' - 8"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 42 1/2"'
$wsPrint.Range("A45").Value = 'LEVEL 2 AREA A Tag: CT-10                                                           TOU: 7'' - 8"                                                                    DEEP STRUT: 14"                                                               Allthread Length: 42 1/2"'
$wsPrint.Range("A46").Value = 'LEVEL 2 AREA A Tag: CT-11                                                           TOU: 10'' - 6 1/2"                                                                    DEEP STRUT: 20"                                                               Allthread Length: 8"'
$wsPrint.Range("A47").Value = 'LEVEL 2 AREA A Tag: CT-12                                                           TOU: 9'' - 9"                                                                    DEEP STRUT: 20"                                                               Allthread Length: 17 1/2"'
$wsPrint.Range("A48").Value = 'LEVEL 2 AREA A Tag: CT-12                                                           TOU: 9'' - 9"                                                                    DEEP STRUT: 20"                                                               Allthread Length: 17 1/2"'
$wsPrint.Range("A49").Value = 'LEVEL 2 AREA A Tag: CT-12                                                           TOU: 9'' - 9"                                                                    DEEP STRUT: 20"                                                               Allthread Length: 17 1/2"'
$wsPrint.Range("A50").Value = 'LEVEL 2 AREA A Tag: CT-12                                                           TOU: 9'' - 9"                                                                    DEEP STRUT: 20"                                                               Allthread Length: 17 1/2"'
$wsPrint.Range("A51").Value = 'LEVEL 2 AREA A Tag: CT-13                                                           TOU: 7'' - 8"                                                                    DEEP STRUT: 22"                                                               Allthread Length: 42 1/2"'
$wsPrint.Range("A52").Value = 'LEVEL 2 AREA A Tag: CT-13                                                           TOU: 7'' - 8"                                                                    BACK TO BACK SHALLOW: 22"                                                               Allthread Length: 42 1/2"'

# --- Reselect the first sheet as active ---
$wb.Worksheets.Item(1).Activate()
